# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" positioned right before "总计",
#    populated with the per-fund holdings table for that quarter.
# 2. Prepend a new row to the "总计" (totals) sheet summarizing the
#    2022-Q1 quarter, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q1" sheet, positioned just before "总计"
# ---------------------------------------------------------------------

# Worksheets.Add() always inserts a brand new sheet before the active
# sheet (index 1), shifting every other sheet's index along the way -
# so look up "总计" only *after* the Add() call, by name rather than a
# previously-captured (index-based) reference.
$newWs = $wb.Worksheets.Add()
$newWs.Name = "2022-Q1"

$zongji = $wb.Worksheets.Item("总计")
$tmpWs = $wb.Worksheets.Item("2022-Q1")
$tmpWs.Move($zongji)

# Worksheet handles are index-bound, so re-resolve both sheets by name
# once more now that Move() has shuffled the sheet order again - $tmpWs
# / $zongji above would otherwise silently point at the wrong sheet (or
# nowhere) for the remainder of the script.
$zongji = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Item("2022-Q1")

# Copy the header-row style (bold + border + centered, style index 2 in
# the original workbook) from the "总计" sheet's header so the new
# sheet matches the look of the other quarterly sheets.
$zongji.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Copy the row-index column style (style index 2) down for rows 2-6.
$zongji.Range("A2").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)

# Columns B-G hold text (fund codes, names and numbers-as-text such as
# "3.00" / "010715" that must keep their formatting / leading zeros),
# so force a text number format before writing them. Column H holds a
# real numeric rank.
$ws.Range("B2:G6").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "519125"
$ws.Range("C2").Value = "浦银安盛消费升级混合A"
$ws.Range("D2").Value = "4.64"
$ws.Range("E2").Value = "89.16"
$ws.Range("F2").Value = "6.48"
$ws.Range("G2").Value = "0.3007"
$ws.Range("H2").Value = 4

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "010715"
$ws.Range("C3").Value = "财通资管消费升级一年持有期混合型证券投资基金A"
$ws.Range("D3").Value = "8.73"
$ws.Range("E3").Value = "73.55"
$ws.Range("F3").Value = "3.00"
$ws.Range("G3").Value = "0.2619"
$ws.Range("H3").Value = 9

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "519176"
$ws.Range("C4").Value = "浦银安盛消费升级混合C"
$ws.Range("D4").Value = "2.33"
$ws.Range("E4").Value = "89.16"
$ws.Range("F4").Value = "6.48"
$ws.Range("G4").Value = "0.1510"
$ws.Range("H4").Value = 4

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "519115"
$ws.Range("C5").Value = "浦银安盛红利精选混合"
$ws.Range("D5").Value = "1.79"
$ws.Range("E5").Value = "88.68"
$ws.Range("F5").Value = "4.67"
$ws.Range("G5").Value = "0.0836"
$ws.Range("H5").Value = 6

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "010716"
$ws.Range("C6").Value = "财通资管消费升级一年持有期混合型证券投资基金C"
$ws.Range("D6").Value = "0.46"
$ws.Range("E6").Value = "73.55"
$ws.Range("F6").Value = "3.00"
$ws.Range("G6").Value = "0.0138"
$ws.Range("H6").Value = 9

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 summary row to the "总计" sheet, pushing
# the existing quarters down by one row.
# ---------------------------------------------------------------------

# Create a new A5 cell (doesn't exist yet) that mirrors the style of
# the existing index-column cells (A2:A4, style index 2).
$zongji.Range("A4").Copy()
$zongji.Range("A5").PasteSpecial(-4122)

# Shift the existing quarter rows down by one (bottom-up, to avoid
# clobbering data not yet moved).
$zongji.Range("B5").Value = "2021-Q1"
$zongji.Range("C5").Value = 4
$zongji.Range("D5").Value = 1.44

$zongji.Range("B4").Value = "2021-Q2"
$zongji.Range("C4").Value = 2
$zongji.Range("D4").Value = 0.01

$zongji.Range("B3").Value = "2021-Q4"
$zongji.Range("C3").Value = 12
$zongji.Range("D3").Value = 1.84

$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 5
$zongji.Range("D2").Value = 0.8100000000000001

# Renumber the index column to match the new row order.
$zongji.Range("A2").Value = 0
$zongji.Range("A3").Value = 1
$zongji.Range("A4").Value = 2
$zongji.Range("A5").Value = 3
